# Fractal.Calc.xlsx edit: rename "index" column to "i" and renumber it starting at 0
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fractal(2)")

# 1. Rename header A1 from "index" to "i" (this also renames the table column
#    "testdata[index]" -> "testdata[i]" automatically because A1 is the table header).
$ws.Range("A1").Value = "i"

# 2. Renumber the index column: old values 1..502 (rows 2..503) become 0..501.
$rng = $ws.Range("A2:A503")
$vals = $rng.Value()
$rows = $vals.GetLength(0)
for ($i = 1; $i -le $rows; $i++) {
    $vals[$i, 1] = $vals[$i, 1] - 1
}
$rng.Value = $vals

# 3. Shrink column A width from 6 to 4 characters.
$ws.Columns.Item(1).ColumnWidth = 3.14

# 4. Move the active selection from I1 to M1.
$ws.Range("M1").Select()
